# Applies the SCD0338 -> SCD0026 rename edit described in the commit
# "Update SCD0026-001 until SCD0026-017 Fix"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/sheet tab from SCD0338 to SCD0026
$ws.Name = "SCD0026"

# Update the TC_ID cells (B2/B3) that carried the old test-case code
$ws.Range("B2").Value = "SCD0026-010"
$ws.Range("B3").Value = "SCD0026-010"

# Move the active selection from N2 to B4
$ws.Range("B4").Select()
